# "9th Stab - Cosmetic Changes"
#
# The MarketBeat rank history sheet gains two more weekly snapshot columns
# ("Jun_17" and "Jun_15"). They are inserted immediately to the left of the
# existing rating columns, pushing the old "Jun_13" / "Jun_10" columns from
# B/C to D/E. The two new columns are seeded with the same "UN" placeholder
# used throughout the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh columns at B, shifting the old B ("Jun_13") and
# C ("Jun_10") columns right to D and E.
$ws.Columns("B:C").Insert()

# Match the sheet's existing explicit column width (stored as 8.0 chars)
# for the two new columns and the newly-shifted-over old "C" column.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14

# New header row values for the two newly inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Seed the new columns' data rows with the same "UN" placeholder used by
# every other rank cell in the table.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
